# Add two new columns (I: I0, J: IF) to Sheet1, matching existing
# column style (header row uses the bold/bordered style copied from H1;
# data rows use the default style, same as the other data columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from the existing
# header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2-44) ----------------------------------------------------
$I0vals = @(4,9,9,6,5,11,8,7,4,1,8,3,1,5,6,5,7,8,7,6,8,3,7,9,7,6,8,6,7,7,9,5,8,6,7,6,7,4,6,8,5,2,4)
$IFvals = @(6,9,9,8,7,11,9,8,6,4,8,6,3,6,8,5,8,8,7,7,8,5,8,9,8,7,8,6,7,7,9,7,8,6,7,6,7,5,7,8,5,3,4)

for ($i = 0; $i -lt $I0vals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0vals[$i]
    $ws.Cells.Item($row, 10).Value = $IFvals[$i]
}
